$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.076113267444826
$ws.Range("D2").Value = 1.06609586305523
$ws.Range("E2").Value = 1.089871793792328
$ws.Range("F2").Value = 1.097128008043839
$ws.Range("I2").Value = 1.046884506009596
$ws.Range("J2").Value = 1.081014349518926
$ws.Range("K2").Value = 1.068807332416748
$ws.Range("L2").Value = 1.092520799621132
$ws.Range("M2").Value = 1.099758531944935
$ws.Range("N2").Value = 1.082549513889321
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.078446746935534
$ws.Range("D3").Value = 1.067845383074219
$ws.Range("E3").Value = 1.092244857109299
$ws.Range("F3").Value = 1.099651626848998
$ws.Range("I3").Value = 1.047532280513025
$ws.Range("J3").Value = 1.083001123662674
$ws.Range("K3").Value = 1.070370048895124
$ws.Range("L3").Value = 1.094710179142661
$ws.Range("M3").Value = 1.102099493563216
$ws.Range("N3").Value = 1.08453910948024
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.079948072950194
$ws.Range("D4").Value = 1.068969645730607
$ws.Range("E4").Value = 1.093772066259355
$ws.Range("F4").Value = 1.101276287499461
$ws.Range("I4").Value = 1.047946267867412
$ws.Range("J4").Value = 1.08427802461042
$ws.Range("K4").Value = 1.071372996554956
$ws.Range("L4").Value = 1.096118205526243
$ws.Range("M4").Value = 1.10360566018826
$ws.Range("N4").Value = 1.085817823773794
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.08057722206147
$ws.Range("D5").Value = 1.069440456496829
$ws.Range("E5").Value = 1.094412161041039
$ws.Range("F5").Value = 1.101957361117595
$ws.Range("I5").Value = 1.048119087644886
$ws.Range("J5").Value = 1.084812799691846
$ws.Range("K5").Value = 1.071792699402334
$ws.Range("L5").Value = 1.096708115236919
$ws.Range("M5").Value = 1.104236843991751
$ws.Range("N5").Value = 1.086353358297176
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.080682742590467
$ws.Range("D6").Value = 1.069519401482045
$ws.Range("E6").Value = 1.094519523167603
$ws.Range("F6").Value = 1.102071604444496
$ws.Range("I6").Value = 1.048148033743125
$ws.Range("J6").Value = 1.084902472731635
$ws.Range("K6").Value = 1.07186305681656
$ws.Range("L6").Value = 1.096807046273321
$ws.Range("M6").Value = 1.104342706133961
$ws.Range("N6").Value = 1.086443158682966
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.079956487498297
$ws.Range("D7").Value = 1.068975943866002
$ws.Range("E7").Value = 1.093780626811221
$ws.Range("F7").Value = 1.101285395568967
$ws.Range("I7").Value = 1.047948581867778
$ws.Range("J7").Value = 1.084285178236852
$ws.Range("K7").Value = 1.071378612206804
$ws.Range("L7").Value = 1.096126095822688
$ws.Range("M7").Value = 1.103614101925331
$ws.Range("N7").Value = 1.085824987559197
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.076903686030892
$ws.Range("D8").Value = 1.06668875750798
$ws.Range("E8").Value = 1.090675533954442
$ws.Range("F8").Value = 1.097982624056763
$ws.Range("I8").Value = 1.047104503979928
$ws.Range("J8").Value = 1.081687611224379
$ws.Range("K8").Value = 1.069337185993989
$ws.Range("L8").Value = 1.093262529290263
$ws.Range("M8").Value = 1.100551479259917
$ws.Range("N8").Value = 1.083223731703621
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.071456257860726
$ws.Range("D9").Value = 1.062597086039532
$ws.Range("E9").Value = 1.085138035894708
$ws.Range("F9").Value = 1.092096920130095
$ws.Range("I9").Value = 1.045576824026608
$ws.Range("J9").Value = 1.07704198518203
$ws.Range("K9").Value = 1.065675273824696
$ws.Range("L9").Value = 1.088148232487491
$ws.Range("M9").Value = 1.095086743215559
$ws.Range("N9").Value = 1.078571508339431
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.067775798023058
$ws.Range("D10").Value = 1.059825691630359
$ws.Range("E10").Value = 1.081398935793902
$ws.Range("F10").Value = 1.088125605500044
$ws.Range("I10").Value = 1.044530219344451
$ws.Range("J10").Value = 1.073896173408475
$ws.Range("K10").Value = 1.063188269233177
$ws.Range("L10").Value = 1.084689803240909
$ws.Range("M10").Value = 1.091394750237903
$ws.Range("N10").Value = 1.075421229152328
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.066169827730754
$ws.Range("D11").Value = 1.058614768210669
$ws.Range("E11").Value = 1.079767909863133
$ws.Range("F11").Value = 1.086393972744472
$ws.Range("I11").Value = 1.044070112379128
$ws.Range("J11").Value = 1.07252182160457
$ws.Range("K11").Value = 1.062100011093959
$ws.Range("L11").Value = 1.083179997535802
$ws.Range("M11").Value = 1.089783788733937
$ws.Range("N11").Value = 1.074044925611221
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.065571388332636
$ws.Range("D12").Value = 1.058163294393495
$ws.Range("E12").Value = 1.079160214399706
$ws.Range("F12").Value = 1.085748895737457
$ws.Range("I12").Value = 1.043898147625207
$ws.Range("J12").Value = 1.07200943971965
$ws.Range("K12").Value = 1.061694030924341
$ws.Range("L12").Value = 1.082617284693836
$ws.Range("M12").Value = 1.089183495912251
$ws.Range("N12").Value = 1.073531816085253
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.065699843074542
$ws.Range("D13").Value = 1.05826021397755
$ws.Range("E13").Value = 1.079290652300226
$ws.Range("F13").Value = 1.085887352687382
$ws.Range("I13").Value = 1.043935082980766
$ws.Range("J13").Value = 1.072119433624894
$ws.Range("K13").Value = 1.061781195087625
$ws.Range("L13").Value = 1.082738075611337
$ws.Range("M13").Value = 1.089312348193535
$ws.Range("N13").Value = 1.073641966194459
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.066120399878714
$ws.Range("D14").Value = 1.058577483849291
$ws.Range("E14").Value = 1.079717715865891
$ws.Range("F14").Value = 1.086340689034494
$ws.Range("I14").Value = 1.044055919476043
$ws.Range("J14").Value = 1.0724795067596
$ws.Range("K14").Value = 1.062066488661669
$ws.Range("L14").Value = 1.083133522695876
$ws.Range("M14").Value = 1.089734207665674
$ws.Range("N14").Value = 1.074002550674317
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.06637926377825
$ws.Range("D15").Value = 1.058772739869448
$ws.Range("E15").Value = 1.079980595523292
$ws.Range("F15").Value = 1.086619754573334
$ws.Range("I15").Value = 1.044130229665047
$ws.Range("J15").Value = 1.072701108085278
$ws.Range("K15").Value = 1.062242033854852
$ws.Range("L15").Value = 1.083376916658104
$ws.Range("M15").Value = 1.089993874256638
$ws.Range("N15").Value = 1.074224466699296
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.067882116467017
$ws.Range("D16").Value = 1.059905823203915
$ws.Range("E16").Value = 1.081506924196889
$ws.Range("F16").Value = 1.08824026938544
$ws.Range("I16").Value = 1.044560607557853
$ws.Range("J16").Value = 1.073987123075525
$ws.Range("K16").Value = 1.06326025002975
$ws.Range("L16").Value = 1.084789740329112
$ws.Range("M16").Value = 1.091501400000085
$ws.Range("N16").Value = 1.075512307978337
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.068821478318813
$ws.Range("D17").Value = 1.060613627917339
$ws.Range("E17").Value = 1.082461102030797
$ws.Range("F17").Value = 1.089253510346945
$ws.Range("I17").Value = 1.04482870492318
$ws.Range("J17").Value = 1.0747905034916
$ws.Range("K17").Value = 1.063895875213434
$ws.Range("L17").Value = 1.085672637958201
$ws.Range("M17").Value = 1.092443693950854
$ws.Range("N17").Value = 1.076316829286752
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.069368208982027
$ws.Range("D18").Value = 1.061025430775725
$ws.Range("E18").Value = 1.083016507006811
$ws.Range("F18").Value = 1.089843360624667
$ws.Range("I18").Value = 1.04498441531088
$ws.Range("J18").Value = 1.075257928829165
$ws.Range("K18").Value = 1.064265530773639
$ws.Range("L18").Value = 1.086186436082603
$ws.Range("M18").Value = 1.092992135110949
$ws.Range("N18").Value = 1.076784918421911
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.069554431158958
$ws.Range("D19").Value = 1.061165668604098
$ws.Range("E19").Value = 1.083205692388965
$ws.Range("F19").Value = 1.090044290175249
$ws.Range("I19").Value = 1.045037396142331
$ws.Range("J19").Value = 1.07541711147809
$ws.Range("K19").Value = 1.064391389565653
$ws.Range("L19").Value = 1.086361429374195
$ws.Range("M19").Value = 1.093178940562251
$ws.Range("N19").Value = 1.076944327128455
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.068720816508579
$ws.Range("D20").Value = 1.060537795868491
$ws.Range("E20").Value = 1.082358847252815
$ws.Range("F20").Value = 1.089144919217461
$ws.Range("I20").Value = 1.044800009677993
$ws.Range("J20").Value = 1.07470443008073
$ws.Range("K20").Value = 1.063827792109157
$ws.Range("L20").Value = 1.085578033986175
$ws.Range("M20").Value = 1.092342717527059
$ws.Range("N20").Value = 1.076230633641766
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.065996609622769
$ws.Range("D21").Value = 1.058484102616846
$ws.Range("E21").Value = 1.079592008128962
$ws.Range("F21").Value = 1.086207244967105
$ws.Range("I21").Value = 1.044020365585843
$ws.Range("J21").Value = 1.072373526712997
$ws.Range("K21").Value = 1.061982525559497
$ws.Range("L21").Value = 1.083017126414571
$ws.Range("M21").Value = 1.089610033790686
$ws.Range("N21").Value = 1.073896420123892
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.064272707451794
$ws.Range("D22").Value = 1.057183103325185
$ws.Range("E22").Value = 1.077841594918188
$ws.Range("F22").Value = 1.084349351658276
$ws.Range("I22").Value = 1.043524025769324
$ws.Range("J22").Value = 1.070897053158369
$ws.Range("K22").Value = 1.060812170181494
$ws.Range("L22").Value = 1.081395937652392
$ws.Range("M22").Value = 1.087880805045963
$ws.Range("N22").Value = 1.072417849807463
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.065187652262914
$ws.Range("D23").Value = 1.057873728721874
$ws.Range("E23").Value = 1.078770565827964
$ws.Range("F23").Value = 1.085335307459726
$ws.Range("I23").Value = 1.043787734613191
$ws.Range("J23").Value = 1.071680815366387
$ws.Range("K23").Value = 1.06143357620775
$ws.Range("L23").Value = 1.082256427216342
$ws.Range("M23").Value = 1.088798573231987
$ws.Range("N23").Value = 1.073202725047715
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.068766304900776
$ws.Range("D24").Value = 1.060572064342205
$ws.Range("E24").Value = 1.08240505533848
$ws.Range("F24").Value = 1.089193990443599
$ws.Range("I24").Value = 1.044812977881865
$ws.Range("J24").Value = 1.074743326571592
$ws.Range("K24").Value = 1.063858559311188
$ws.Range("L24").Value = 1.085620785095754
$ws.Range("M24").Value = 1.092388348078814
$ws.Range("N24").Value = 1.076269585370106
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.072872924470525
$ws.Range("D25").Value = 1.063662392273617
$ws.Range("E25").Value = 1.086577740736948
$ws.Range("F25").Value = 1.093626645927543
$ws.Range("I25").Value = 1.045976649313953
$ws.Range("J25").Value = 1.078251374186791
$ws.Range("K25").Value = 1.066629857701387
$ws.Range("L25").Value = 1.089478798721514
$ws.Range("M25").Value = 1.096507883856555
$ws.Range("N25").Value = 1.079782614815298
